$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header for the new "비추천수" (dislikes) column, styled like the other headers
$ws.Range("F1").Value = "비추천수"
$ws.Range("F1").Style = $ws.Range("E1").Style

# Dislike counts for rows 2 through 114 (as text, matching the existing "추천수" column format)
$values = @("0", "0", "1", "1", "0", "0", "0", "0", "0", "0", "0", "1", "1", "0", "1", "0", "0", "0", "0", "0", "4", "0", "5", "0", "0", "2", "0", "0", "0", "3", "0", "1", "0", "0", "0", "0", "2", "1", "0", "0", "2", "8", "0", "0", "3", "12", "3", "0", "0", "0", "0", "4", "0", "0", "0", "4", "1", "1", "0", "0", "3", "1", "1", "0", "0", "4", "0", "0", "1", "0", "0", "0", "3", "0", "3", "0", "0", "0", "0", "0", "1", "0", "0", "0", "1", "0", "0", "0", "1", "0", "0", "0", "0", "1", "0", "2", "0", "1", "0", "0", "0", "1", "0", "5", "0", "0", "1", "4", "0", "3", "15", "9", "2")

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $values[$i]
}

Write-Host "Done adding dislikes column"
